$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.155.49'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.73%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.831.46'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.19%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.21'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.09%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4647'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.91%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2691'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -6.95%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06260'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.42%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.841.45'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.19%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07383'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.09%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.03'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.75%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.892'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.44%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '83.25'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -5.45%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6195'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -7.34%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.084.67'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.84%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9997'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.23'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.60%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007274'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.16%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -6.59%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.070.18'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.99%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.839'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -8.39%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("B24").Style = "Normal"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("B24").Style = "Normal"

$ws.Range("C24").Style = "Normal"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("C24").Style = "Normal"

$ws.Range("D24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.867'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.19%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("B25").Style = "Normal"
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("B25").Style = "Normal"

$ws.Range("C25").Style = "Normal"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("C25").Style = "Normal"

$ws.Range("D25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.113'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.01%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.35'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.10%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("B27").Style = "Normal"
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("B27").Style = "Normal"

$ws.Range("C27").Style = "Normal"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("C27").Style = "Normal"

$ws.Range("D27").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.65'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -5.91%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").Style = "Normal"
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("B28").Style = "Normal"

$ws.Range("C28").Style = "Normal"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("C28").Style = "Normal"

$ws.Range("D28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.842'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.98%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("B29").Style = "Normal"
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'Stellar'
$ws.Range("B29").Style = "Normal"

$ws.Range("C29").Style = "Normal"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("C29").Style = "Normal"

$ws.Range("D29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1007'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.72%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("B30").Style = "Normal"
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("B30").Style = "Normal"

$ws.Range("C30").Style = "Normal"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("C30").Style = "Normal"

$ws.Range("D30").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.367'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.44%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("B31").Style = "Normal"
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("B31").Style = "Normal"

$ws.Range("C31").Style = "Normal"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("C31").Style = "Normal"

$ws.Range("D31").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.043'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -6.99%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").Style = "Normal"
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("B32").Style = "Normal"

$ws.Range("C32").Style = "Normal"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("C32").Style = "Normal"

$ws.Range("D32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.766'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -6.69%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").Style = "Normal"
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'Hedera'
$ws.Range("B33").Style = "Normal"

$ws.Range("C33").Style = "Normal"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("C33").Style = "Normal"

$ws.Range("D33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04776'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -5.93%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").Style = "Normal"
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("B34").Style = "Normal"

$ws.Range("C34").Style = "Normal"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("C34").Style = "Normal"

$ws.Range("D34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.124'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -7.26%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("B35").Style = "Normal"
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("B35").Style = "Normal"

$ws.Range("C35").Style = "Normal"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C35").Style = "Normal"

$ws.Range("D35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6983'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -7.06%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("B36").Style = "Normal"
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("B36").Style = "Normal"

$ws.Range("C36").Style = "Normal"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("C36").Style = "Normal"

$ws.Range("D36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.687'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.01%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Style = "Normal"
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'VeChain'
$ws.Range("B37").Style = "Normal"

$ws.Range("C37").Style = "Normal"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("C37").Style = "Normal"

$ws.Range("D37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01804'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.30%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Style = "Normal"
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'MXToken'
$ws.Range("B38").Style = "Normal"

$ws.Range("C38").Style = "Normal"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("C38").Style = "Normal"

$ws.Range("D38").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.609'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.44%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").Style = "Normal"
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("B39").Style = "Normal"

$ws.Range("C39").Style = "Normal"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("C39").Style = "Normal"

$ws.Range("D39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8938'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.01%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Style = "Normal"
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("B40").Style = "Normal"

$ws.Range("C40").Style = "Normal"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("C40").Style = "Normal"

$ws.Range("D40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.922'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -7.04%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Style = "Normal"
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("B41").Style = "Normal"

$ws.Range("C41").Style = "Normal"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("C41").Style = "Normal"

$ws.Range("D41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.000'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Style = "Normal"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'Quant'
$ws.Range("B42").Style = "Normal"

$ws.Range("C42").Style = "Normal"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("C42").Style = "Normal"

$ws.Range("D42").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '102.81'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.85%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Style = "Normal"
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("B43").Style = "Normal"

$ws.Range("C43").Style = "Normal"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C43").Style = "Normal"

$ws.Range("D43").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.489'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.64%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").Style = "Normal"
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("B44").Style = "Normal"

$ws.Range("C44").Style = "Normal"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("C44").Style = "Normal"

$ws.Range("D44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3991'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -7.15%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Style = "Normal"
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Aptos'
$ws.Range("B45").Style = "Normal"

$ws.Range("C45").Style = "Normal"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C45").Style = "Normal"

$ws.Range("D45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.956'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -6.32%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Style = "Normal"
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Algorand'
$ws.Range("B46").Style = "Normal"

$ws.Range("C46").Style = "Normal"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("C46").Style = "Normal"

$ws.Range("D46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1190'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -7.05%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Style = "Normal"
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Aave'
$ws.Range("B47").Style = "Normal"

$ws.Range("C47").Style = "Normal"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("C47").Style = "Normal"

$ws.Range("D47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '59.53'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -7.41%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Style = "Normal"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("B48").Style = "Normal"

$ws.Range("C48").Style = "Normal"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("C48").Style = "Normal"

$ws.Range("D48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.412'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -6.93%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Style = "Normal"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Cronos'
$ws.Range("B49").Style = "Normal"

$ws.Range("C49").Style = "Normal"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("C49").Style = "Normal"

$ws.Range("D49").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05523'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.68%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Style = "Normal"
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Elrond'
$ws.Range("B50").Style = "Normal"

$ws.Range("C50").Style = "Normal"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("C50").Style = "Normal"

$ws.Range("D50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '32.51'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.52%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Style = "Normal"
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("B51").Style = "Normal"

$ws.Range("C51").Style = "Normal"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("C51").Style = "Normal"

$ws.Range("D51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.358'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -8.91%  '
$ws.Range("E51").Style = "Normal"
